$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 150.33333
$ws.Range("J9").Value = 134.33333
$ws.Range("L9").Value = 134.33333
$ws.Range("N9").Value = -472.33333
$ws.Range("H17").Value = 688.35486
$ws.Range("J17").Value = 704.63336
$ws.Range("L17").Value = 2113.90008
$ws.Range("N17").Value = -2449.90008
$ws.Range("H41").Value = 1021.2857
$ws.Range("I41").Value = 2199
$ws.Range("J41").Value = 550.2
$ws.Range("K41").Value = 2199
$ws.Range("L41").Value = 550.2
$ws.Range("M41").Value = -1759
$ws.Range("N41").Value = -1430.2
$ws.Range("H43").Value = 3180.6667
$ws.Range("I43").Value = 2811
$ws.Range("J43").Value = 3365.5
$ws.Range("K43").Value = 2811
$ws.Range("L43").Value = 3365.5
$ws.Range("M43").Value = -2742
$ws.Range("N43").Value = -3503.5
$ws.Range("H70").Value = 1290.5454
$ws.Range("I70").Value = 966.2222
$ws.Range("K70").Value = 2898.6666
$ws.Range("M70").Value = -2628.6666
$ws.Range("H73").Value = 1290.5454
$ws.Range("I73").Value = 966.2222
$ws.Range("K73").Value = 2898.6666
$ws.Range("M73").Value = -1962.6666
$ws.Range("H76").Value = 4799.5
$ws.Range("J76").Value = 4856.857
$ws.Range("L76").Value = 4856.857
$ws.Range("N76").Value = -5486.857
$ws.Range("H79").Value = 4799.5
$ws.Range("J79").Value = 4856.857
$ws.Range("L79").Value = 4856.857
$ws.Range("N79").Value = -7040.857
$ws.Range("H80").Value = 2785.55
$ws.Range("I80").Value = 339.7143
$ws.Range("K80").Value = 1019.1429
$ws.Range("M80").Value = -21.14289999999994
$ws.Range("H83").Value = 2785.55
$ws.Range("I83").Value = 339.7143
$ws.Range("K83").Value = 3057.4287
$ws.Range("M83").Value = 1934.5713
$ws.Range("H86").Value = 1920
$ws.Range("I86").Value = 1982.2858
$ws.Range("J86").Value = 1484
$ws.Range("K86").Value = 1982.2858
$ws.Range("L86").Value = 1484
$ws.Range("M86").Value = -859.2858000000001
$ws.Range("N86").Value = -3730
$ws.Range("H88").Value = 4832.5
$ws.Range("J88").Value = 4599.2
$ws.Range("L88").Value = 4599.2
$ws.Range("N88").Value = -5411.2
$ws.Range("H89").Value = 1920
$ws.Range("I89").Value = 1982.2858
$ws.Range("J89").Value = 1484
$ws.Range("K89").Value = 9911.429
$ws.Range("L89").Value = 7420
$ws.Range("M89").Value = -4295.429
$ws.Range("N89").Value = -18652
$ws.Range("H91").Value = 4832.5
$ws.Range("J91").Value = 4599.2
$ws.Range("L91").Value = 4599.2
$ws.Range("N91").Value = -7407.2
$ws.Range("H98").Value = 1932.9375
$ws.Range("I98").Value = 2006
$ws.Range("J98").Value = 1772.2
$ws.Range("K98").Value = 2006
$ws.Range("L98").Value = 1772.2
$ws.Range("M98").Value = -508
$ws.Range("N98").Value = -4768.2
$ws.Range("H100").Value = 4111.5
$ws.Range("I100").Value = 3478.6
$ws.Range("J100").Value = 5166.3335
$ws.Range("K100").Value = 3478.6
$ws.Range("L100").Value = 5166.3335
$ws.Range("M100").Value = -2937.6
$ws.Range("N100").Value = -6248.3335
$ws.Range("H106").Value = 1799
$ws.Range("I106").Value = 1799
$ws.Range("K106").Value = 1799
$ws.Range("M106").Value = -1168
$ws.Range("H107").Value = 2997.4
$ws.Range("I107").Value = 2997.4
$ws.Range("K107").Value = 2997.4
$ws.Range("M107").Value = -1077.4
$ws.Range("H108").Value = 50000
$ws.Range("I108").Value = 50000
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 50000
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("M108").Value = -46160
$ws.Range("H113").Value = 2618.75
$ws.Range("I113").Value = 2618.75
$ws.Range("K113").Value = 2618.75
$ws.Range("M113").Value = 635.25
$ws.Range("H114").Value = 68073.664
$ws.Range("I114").Value = 67860.5
$ws.Range("J114").Value = 68500
$ws.Range("K114").Value = 67860.5
$ws.Range("L114").Value = 68500
$ws.Range("M114").Value = -63521.5
$ws.Range("N114").Value = -77178
$ws.Range("H116").Value = 44093.406
$ws.Range("I116").Value = 6197.1333
$ws.Range("K116").Value = 6197.1333
$ws.Range("M116").Value = -2755.1333
$ws.Range("H122").Value = 1932.9375
$ws.Range("I122").Value = 2006
$ws.Range("J122").Value = 1772.2
$ws.Range("K122").Value = 6018
$ws.Range("L122").Value = 5316.6
$ws.Range("M122").Value = -3568
$ws.Range("N122").Value = -10216.6
$ws.Range("H126").Value = 120000
$ws.Range("J126").Value = 120000
$ws.Range("L126").Value = 120000
$ws.Range("N126").Value = -129880
$ws.Range("H127").Value = 1950
$ws.Range("I127").Value = 1950
$ws.Range("K127").Value = 5850
$ws.Range("M127").Value = -890
$ws.Range("H131").Value = 100000000
$ws.Range("I131").Value = 100000000
$ws.Range("K131").Value = 300000000
$ws.Range("M131").Value = -299994960
$ws.Range("H132").Value = 3129.1765
$ws.Range("I132").Value = 3092.8
$ws.Range("J132").Value = 3402
$ws.Range("K132").Value = 9278.400000000001
$ws.Range("L132").Value = 10206
$ws.Range("M132").Value = -6748.400000000001
$ws.Range("N132").Value = -15266
$ws.Range("H137").Value = 2465.4333
$ws.Range("I137").Value = 1578.6666
$ws.Range("J137").Value = 4534.5557
$ws.Range("K137").Value = 4735.9998
$ws.Range("L137").Value = 13603.6671
$ws.Range("M137").Value = -2185.9998
$ws.Range("N137").Value = -18703.6671
$ws.Range("H138").Value = 1981
$ws.Range("I138").Value = 1842.9546
$ws.Range("J138").Value = 3499.5
$ws.Range("K138").Value = 5528.8638
$ws.Range("L138").Value = 10498.5
$ws.Range("M138").Value = -388.8638000000001
$ws.Range("N138").Value = -20778.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 22742.75
$ws.Range("I43").Value = 20198
$ws.Range("J43").Value = 30377
$ws.Range("K43").Value = 20198
$ws.Range("L43").Value = 30377
$ws.Range("M43").Value = -19885
$ws.Range("N43").Value = -31003
$ws.Range("H45").Value = 2186
$ws.Range("I45").Value = 2183
$ws.Range("K45").Value = 2183
$ws.Range("M45").Value = -1806
$ws.Range("H61").Value = 3327.025
$ws.Range("I61").Value = 3312.182
$ws.Range("J61").Value = 3397
$ws.Range("K61").Value = 3312.182
$ws.Range("L61").Value = 3397
$ws.Range("M61").Value = -3100.182
$ws.Range("N61").Value = -3821
$ws.Range("H63").Value = 1983.9445
$ws.Range("I63").Value = 2081.3333
$ws.Range("J63").Value = 1935.25
$ws.Range("K63").Value = 2081.3333
$ws.Range("L63").Value = 1935.25
$ws.Range("M63").Value = -1395.3333
$ws.Range("N63").Value = -3307.25
$ws.Range("H66").Value = 1983.9445
$ws.Range("I66").Value = 2081.3333
$ws.Range("J66").Value = 1935.25
$ws.Range("K66").Value = 10406.6665
$ws.Range("L66").Value = 9676.25
$ws.Range("M66").Value = -6974.666499999999
$ws.Range("N66").Value = -16540.25
$ws.Range("H74").Value = 2895117.5
$ws.Range("I74").Value = 1853570.6
$ws.Range("J74").Value = 4631029
$ws.Range("K74").Value = 1853570.6
$ws.Range("L74").Value = 4631029
$ws.Range("M74").Value = -1852696.6
$ws.Range("N74").Value = -4632777
$ws.Range("H76").Value = 48333
$ws.Range("J76").Value = 48333
$ws.Range("L76").Value = 48333
$ws.Range("N76").Value = -49009
$ws.Range("H77").Value = 2895117.5
$ws.Range("I77").Value = 1853570.6
$ws.Range("J77").Value = 4631029
$ws.Range("K77").Value = 9267853
$ws.Range("L77").Value = 23155145
$ws.Range("M77").Value = -9263485
$ws.Range("N77").Value = -23163881
$ws.Range("H79").Value = 48333
$ws.Range("J79").Value = 48333
$ws.Range("L79").Value = 48333
$ws.Range("N79").Value = -50673
$ws.Range("H88").Value = 6399.8
$ws.Range("J88").Value = 5824.75
$ws.Range("L88").Value = 5824.75
$ws.Range("N88").Value = -6636.75
$ws.Range("H91").Value = 6399.8
$ws.Range("J91").Value = 5824.75
$ws.Range("L91").Value = 5824.75
$ws.Range("N91").Value = -8632.75
$ws.Range("H97").Value = 642.7778
$ws.Range("I97").Value = 915.6
$ws.Range("J97").Value = 301.75
$ws.Range("K97").Value = 915.6
$ws.Range("L97").Value = 301.75
$ws.Range("M97").Value = -419.6
$ws.Range("N97").Value = -1293.75
$ws.Range("H122").Value = 2053.5
$ws.Range("I122").Value = 1983.2
$ws.Range("K122").Value = 5949.6
$ws.Range("M122").Value = -3499.6
$ws.Range("H132").Value = 11630952
$ws.Range("I132").Value = 2767.9714
$ws.Range("J132").Value = 62504256
$ws.Range("K132").Value = 8303.914199999999
$ws.Range("L132").Value = 187512768
$ws.Range("M132").Value = -5773.914199999999
$ws.Range("N132").Value = -187517828
$ws.Range("H136").Value = 3327.025
$ws.Range("I136").Value = 3312.182
$ws.Range("J136").Value = 3397
$ws.Range("K136").Value = 9936.545999999998
$ws.Range("L136").Value = 10191
$ws.Range("M136").Value = -7386.545999999998
$ws.Range("N136").Value = -15291
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1328.2858
$ws.Range("I64").Value = 900
$ws.Range("K64").Value = 900
$ws.Range("M64").Value = -675
$ws.Range("H67").Value = 1328.2858
$ws.Range("I67").Value = 900
$ws.Range("K67").Value = 900
$ws.Range("M67").Value = -120
$ws.Range("H86").Value = 3862
$ws.Range("I86").Value = 2003.5
$ws.Range("K86").Value = 2003.5
$ws.Range("M86").Value = -880.5
$ws.Range("H89").Value = 3862
$ws.Range("I89").Value = 2003.5
$ws.Range("K89").Value = 10017.5
$ws.Range("M89").Value = -4401.5
$ws.Range("H94").Value = 244.14285
$ws.Range("I94").Value = 170.875
$ws.Range("J94").Value = 341.83334
$ws.Range("K94").Value = 170.875
$ws.Range("L94").Value = 341.83334
$ws.Range("M94").Value = 280.125
$ws.Range("N94").Value = -1243.83334
$ws.Range("H126").Value = 90590
$ws.Range("J126").Value = 90590
$ws.Range("L126").Value = 90590
$ws.Range("N126").Value = -100470
$ws.Range("H134").Value = 27781050
$ws.Range("I134").Value = 13892243
$ws.Range("J134").Value = 111113896
$ws.Range("K134").Value = 41676729
$ws.Range("L134").Value = 333341688
$ws.Range("M134").Value = -41674194
$ws.Range("N134").Value = -333346758
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 25000476
$ws.Range("I4").Value = 950
$ws.Range("K4").Value = 950
$ws.Range("M4").Value = -838
$ws.Range("H16").Value = 2996
$ws.Range("I16").Value = 2996
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2996
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2709
$ws.Range("N16").ClearContents()
$ws.Range("H58").Value = 2629.389
$ws.Range("I58").Value = 2509.0833
$ws.Range("K58").Value = 2509.0833
$ws.Range("M58").Value = -2306.0833
$ws.Range("H99").Value = 4024.6667
$ws.Range("I99").Value = 3328.7
$ws.Range("K99").Value = 3328.7
$ws.Range("M99").Value = -1830.7
$ws.Range("H107").Value = 2389.0588
$ws.Range("J107").Value = 3680.25
$ws.Range("L107").Value = 3680.25
$ws.Range("N107").Value = -7520.25
$ws.Range("H113").Value = 2996
$ws.Range("I113").Value = 2996
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2996
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -826
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 4024.6667
$ws.Range("I126").Value = 3328.7
$ws.Range("K126").Value = 9986.099999999999
$ws.Range("M126").Value = -7516.099999999999
$ws.Range("H132").Value = 4648.9375
$ws.Range("I132").Value = 4670.9287
$ws.Range("K132").Value = 14012.7861
$ws.Range("M132").Value = -11482.7861
$ws.Range("H134").Value = 5558838
$ws.Range("I134").Value = 2708.818
$ws.Range("K134").Value = 8126.454000000001
$ws.Range("M134").Value = -5591.454000000001
$ws.Range("H136").Value = 2629.389
$ws.Range("I136").Value = 2509.0833
$ws.Range("K136").Value = 7527.249899999999
$ws.Range("M136").Value = -4977.249899999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 143695.14
$ws.Range("I4").Value = 126.6
$ws.Range("J4").Value = 223455.44
$ws.Range("K4").Value = 379.8
$ws.Range("L4").Value = 670366.3200000001
$ws.Range("M4").Value = -267.8
$ws.Range("N4").Value = -670590.3200000001
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H23").Value = 1452
$ws.Range("I23").Value = 2996
$ws.Range("J23").Value = 569.7143
$ws.Range("K23").Value = 8988
$ws.Range("L23").Value = 1709.1429
$ws.Range("M23").Value = -8753
$ws.Range("N23").Value = -2179.1429
$ws.Range("H109").Value = 3392.6
$ws.Range("I109").Value = 3258.4443
$ws.Range("J109").Value = 4600
$ws.Range("K109").Value = 9775.332900000001
$ws.Range("L109").Value = 13800
$ws.Range("M109").Value = -8735.332900000001
$ws.Range("N109").Value = -15880
$ws.Range("H121").Value = 1910963.1
$ws.Range("I121").Value = 144344.42
$ws.Range("J121").Value = 5002546
$ws.Range("K121").Value = 433033.26
$ws.Range("L121").Value = 15007638
$ws.Range("M121").Value = -431723.26
$ws.Range("N121").Value = -15010258
$ws.Range("H124").Value = 1189
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("H130").Value = 5028.778
$ws.Range("J130").Value = 5028.778
$ws.Range("L130").Value = 15086.334
$ws.Range("N130").Value = -25126.334
$ws.Range("H131").Value = 467633.78
$ws.Range("J131").Value = 506511.6
$ws.Range("L131").Value = 1519534.8
$ws.Range("N131").Value = -1529614.8
$ws.Range("H132").Value = 933
$ws.Range("J132").Value = 2500
$ws.Range("L132").Value = 22500
$ws.Range("N132").Value = -27560
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2787.625
$ws.Range("I80").Value = 2780.4
$ws.Range("J80").Value = 2799.6667
$ws.Range("K80").Value = 2780.4
$ws.Range("L80").Value = 2799.6667
$ws.Range("M80").Value = -1782.4
$ws.Range("N80").Value = -4795.6667
$ws.Range("H83").Value = 2787.625
$ws.Range("I83").Value = 2780.4
$ws.Range("J83").Value = 2799.6667
$ws.Range("K83").Value = 13902
$ws.Range("L83").Value = 13998.3335
$ws.Range("M83").Value = -8910
$ws.Range("N83").Value = -23982.3335
$ws.Range("H97").Value = 571.9167
$ws.Range("I97").Value = 625
$ws.Range("K97").Value = 625
$ws.Range("M97").Value = -129
$ws.Range("H107").Value = 759.5333000000001
$ws.Range("I107").Value = 345.4
$ws.Range("K107").Value = 345.4
$ws.Range("M107").Value = 1574.6
$ws.Range("H122").Value = 1944.9524
$ws.Range("I122").Value = 2203.7407
$ws.Range("J122").Value = 1479.1333
$ws.Range("K122").Value = 6611.222099999999
$ws.Range("L122").Value = 4437.3999
$ws.Range("M122").Value = -4161.222099999999
$ws.Range("N122").Value = -9337.3999
$ws.Range("H126").Value = 8663.75
$ws.Range("I126").Value = 8663.75
$ws.Range("K126").Value = 25991.25
$ws.Range("M126").Value = -23521.25
$ws.Range("H132").Value = 2270.75
$ws.Range("I132").Value = 2216.8386
$ws.Range("J132").Value = 2399.3076
$ws.Range("K132").Value = 6650.5158
$ws.Range("L132").Value = 7197.9228
$ws.Range("M132").Value = -4120.5158
$ws.Range("N132").Value = -12257.9228
$ws.Range("H134").Value = 18000
$ws.Range("J134").Value = 18000
$ws.Range("L134").Value = 54000
$ws.Range("N134").Value = -59070
$ws.Range("H135").Value = 90779
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 90779
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 90779
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -100919
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 81081.25
$ws.Range("J139").Value = 72162.5
$ws.Range("L139").Value = 72162.5
$ws.Range("N139").Value = -82442.5
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2865.2856
$ws.Range("I7").Value = 2058.4
$ws.Range("K7").Value = 2058.4
$ws.Range("M7").Value = -1946.4
$ws.Range("H22").Value = 855
$ws.Range("I22").Value = 439.8
$ws.Range("J22").Value = 1374
$ws.Range("K22").Value = 439.8
$ws.Range("L22").Value = 1374
$ws.Range("M22").Value = -144.8
$ws.Range("N22").Value = -1964
$ws.Range("H27").Value = 855
$ws.Range("I27").Value = 439.8
$ws.Range("J27").Value = 1374
$ws.Range("K27").Value = 439.8
$ws.Range("L27").Value = 1374
$ws.Range("M27").Value = -332.8
$ws.Range("N27").Value = -1588
$ws.Range("H61").Value = 2770.3462
$ws.Range("I61").Value = 2423.5789
$ws.Range("J61").Value = 3711.5715
$ws.Range("K61").Value = 2423.5789
$ws.Range("L61").Value = 3711.5715
$ws.Range("M61").Value = -2221.5789
$ws.Range("N61").Value = -4115.5715
$ws.Range("H68").Value = 2392.6428
$ws.Range("I68").Value = 2169.9
$ws.Range("K68").Value = 2169.9
$ws.Range("M68").Value = -1420.9
$ws.Range("H71").Value = 2392.6428
$ws.Range("I71").Value = 2169.9
$ws.Range("K71").Value = 10849.5
$ws.Range("M71").Value = -7105.5
$ws.Range("H82").Value = 3985
$ws.Range("I82").Value = 3389.111
$ws.Range("K82").Value = 3389.111
$ws.Range("M82").Value = -3028.111
$ws.Range("H85").Value = 3985
$ws.Range("I85").Value = 3389.111
$ws.Range("K85").Value = 3389.111
$ws.Range("M85").Value = -2141.111
$ws.Range("H100").Value = 4907.5713
$ws.Range("I100").Value = 5151.3335
$ws.Range("K100").Value = 5151.3335
$ws.Range("M100").Value = -4610.3335
$ws.Range("H113").Value = 2770.3462
$ws.Range("I113").Value = 2423.5789
$ws.Range("J113").Value = 3711.5715
$ws.Range("K113").Value = 2423.5789
$ws.Range("L113").Value = 3711.5715
$ws.Range("M113").Value = -253.5789
$ws.Range("N113").Value = -8051.5715
$ws.Range("H122").Value = 3535.0588
$ws.Range("I122").Value = 3325
$ws.Range("K122").Value = 9975
$ws.Range("M122").Value = -7525
$ws.Range("H126").Value = 2865.2856
$ws.Range("I126").Value = 2058.4
$ws.Range("K126").Value = 6175.200000000001
$ws.Range("M126").Value = -3705.200000000001
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 29416194
$ws.Range("I136").Value = 3225.3044
$ws.Range("K136").Value = 9675.913199999999
$ws.Range("M136").Value = -7125.913199999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 11005.5
$ws.Range("J30").Value = 20010
$ws.Range("L30").Value = 20010
$ws.Range("N30").Value = -20224
$ws.Range("H41").Value = 11959.875
$ws.Range("J41").Value = 12203
$ws.Range("L41").Value = 12203
$ws.Range("N41").Value = -12983
$ws.Range("H62").Value = 2601.2727
$ws.Range("I62").Value = 2002.1428
$ws.Range("K62").Value = 2002.1428
$ws.Range("M62").Value = -1378.1428
$ws.Range("H65").Value = 2601.2727
$ws.Range("I65").Value = 2002.1428
$ws.Range("K65").Value = 10010.714
$ws.Range("M65").Value = -6890.714
$ws.Range("H80").Value = 50000
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996
$ws.Range("H81").Value = 8079.4443
$ws.Range("I81").Value = 9359.923000000001
$ws.Range("K81").Value = 18719.846
$ws.Range("M81").Value = -17658.846
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 50000
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984
$ws.Range("H84").Value = 8079.4443
$ws.Range("I84").Value = 9359.923000000001
$ws.Range("K84").Value = 93599.23000000001
$ws.Range("M84").Value = -88295.23000000001
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H113").Value = 2771.8823
$ws.Range("I113").Value = 406.7143
$ws.Range("J113").Value = 4427.5
$ws.Range("K113").Value = 1220.1429
$ws.Range("L113").Value = 13282.5
$ws.Range("M113").Value = 949.8571000000002
$ws.Range("N113").Value = -17622.5
$ws.Range("H122").Value = 9998.5
$ws.Range("I122").Value = 9998.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 29995.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -27545.5
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 2428.2856
$ws.Range("I126").Value = 1559.7
$ws.Range("J126").Value = 4599.75
$ws.Range("K126").Value = 4679.1
$ws.Range("L126").Value = 13799.25
$ws.Range("M126").Value = -2209.1
$ws.Range("N126").Value = -18739.25
$ws.Range("H132").Value = 1392.4117
$ws.Range("I132").Value = 1107.5161
$ws.Range("K132").Value = 3322.5483
$ws.Range("M132").Value = -792.5483000000004
$ws.Range("H136").Value = 1081
$ws.Range("I136").Value = 997.2
$ws.Range("K136").Value = 2991.6
$ws.Range("M136").Value = -441.6000000000004
